# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
#
# The match rows below got their id/result/odds columns (B, F:AC) rotated
# among themselves (the row's A/C/D/E "slot" columns — index, Div, Div
# Original Name, Date — stayed put). Each group rotates so that row N
# receives what used to be row (N+1)'s B..AC content, with the last row
# in the group wrapping around to the first.
#
# NOTE: this COM host does not give PowerShell functions their own
# variable scope (no dot-sourcing isolation), so loop counters defined
# inside a helper function clobber a caller's loop counter of the same
# name. To stay safe, everything below is written inline with uniquely
# named loop variables instead of relying on nested helper functions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B is index 2; columns F..AC are indexes 6..29.
$colIndexes = @(2) + (6..29)

$groups = @(
    @(404, 405, 406, 407),
    @(417, 418, 419, 420, 421),
    @(423, 424)
)

foreach ($group in $groups) {
    # Snapshot every row's current B/F:AC values first (read-before-write,
    # since writes would otherwise clobber data later rows still need).
    $groupSnapshot = @{}
    foreach ($snapRow in $group) {
        $rowVals = @()
        foreach ($snapCol in $colIndexes) {
            $rowVals += ,$ws.Cells.Item($snapRow, $snapCol).Value2
        }
        $groupSnapshot[$snapRow] = $rowVals
    }

    $groupSize = $group.Count
    for ($gi = 0; $gi -lt $groupSize; $gi++) {
        $targetRow = $group[$gi]
        $sourceRow = $group[($gi + 1) % $groupSize]
        $sourceVals = $groupSnapshot[$sourceRow]
        for ($gj = 0; $gj -lt $colIndexes.Count; $gj++) {
            $ws.Cells.Item($targetRow, $colIndexes[$gj]).Value2 = $sourceVals[$gj]
        }
    }
}
